# New weekly price record for "Haba" (Vega Modelo de Temuco) is inserted
# at row 12, pushing the existing rows 12-91 down to 13-92 (the sheet's
# dimension grows from A1:R91 to A1:R92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12:91 down to 13:92, leaving a blank row 12 for the new record.
$ws.Rows("12:12").Insert()

# Populate the new row 12 with this week's data.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 45169
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 100112026
$ws.Range("G12").Value = "Haba"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 13000
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 520
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
